$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing content and shift content down, leaving blank rows in between.
$ws.Range("A1").Value = "Contention 1: BUSH DID 9/11"
$ws.Range("A3").Value = "Subpoint a: JET FUEL CAN'T MELT STEEL BEAMS"
$ws.Range("A5").Value = "Evidence (your mother): JET FUEL CANT MELT STEEL BEAMS"
$ws.Range("A7").Value = "Subpoint b: open your eyes sheeple"
$ws.Range("A9").Value = "Evidence (trump): open"
$ws.Range("A11").Value = "Contention 2: something really important"
$ws.Range("A13").Value = "Subpoint a: jdafn"
$ws.Range("A15").Value = "Subpoint b: jasdfn"

# Touch the interleaved blank rows so they materialize in the sheet.
$ws.Rows.Item(2).OutlineLevel = 0
$ws.Rows.Item(4).OutlineLevel = 0
$ws.Rows.Item(6).OutlineLevel = 0
$ws.Rows.Item(8).OutlineLevel = 0
$ws.Rows.Item(10).OutlineLevel = 0
$ws.Rows.Item(12).OutlineLevel = 0
$ws.Rows.Item(14).OutlineLevel = 0

# Widen column A
$ws.Columns.Item(1).ColumnWidth = 53.1666667
